$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 keeps its textual value "HK_R_acc_SD" (shared string table gains a
# duplicate entry inserted before the existing one during regeneration,
# but the visible cell content is unchanged).
$ws.Range("A1").Value = "HK_R_acc_SD"

# Updated (recomputed) statistics for A2:A49
$ws.Range("A2").Value = 99.360146252285205
$ws.Range("A3").Value = 99.360146252285205
$ws.Range("A4").Value = 99.360146252285205
$ws.Range("A5").Value = 99.360146252285205
$ws.Range("A6").Value = 99.360146252285205
$ws.Range("A7").Value = 99.360146252285205
$ws.Range("A8").Value = 99.268738574040214
$ws.Range("A9").Value = 99.268738574040214
$ws.Range("A10").Value = 99.268738574040214
$ws.Range("A11").Value = 99.268738574040214
$ws.Range("A12").Value = 99.268738574040214
$ws.Range("A13").Value = 99.268738574040214
$ws.Range("A14").Value = 99.268738574040214
$ws.Range("A15").Value = 99.268738574040214
$ws.Range("A16").Value = 99.268738574040214
$ws.Range("A17").Value = 99.268738574040214
$ws.Range("A18").Value = 99.268738574040214
$ws.Range("A19").Value = 99.268738574040214
$ws.Range("A20").Value = 99.634369287020107
$ws.Range("A21").Value = 99.634369287020107
$ws.Range("A22").Value = 99.634369287020107
$ws.Range("A23").Value = 99.360146252285205
$ws.Range("A24").Value = 99.360146252285205
$ws.Range("A25").Value = 99.360146252285205
$ws.Range("A26").Value = 99.314442413162709
$ws.Range("A27").Value = 99.314442413162709
$ws.Range("A28").Value = 99.314442413162709
$ws.Range("A29").Value = 99.268738574040214
$ws.Range("A30").Value = 99.268738574040214
$ws.Range("A31").Value = 99.268738574040214
$ws.Range("A32").Value = 99.405850091407672
$ws.Range("A33").Value = 99.405850091407672
$ws.Range("A34").Value = 99.360146252285205
$ws.Range("A35").Value = 99.360146252285205
$ws.Range("A36").Value = 99.360146252285205
$ws.Range("A37").Value = 99.360146252285205
$ws.Range("A38").Value = 99.360146252285205
$ws.Range("A39").Value = 99.360146252285205
$ws.Range("A40").Value = 99.405850091407672
$ws.Range("A41").Value = 99.268738574040214
$ws.Range("A42").Value = 99.268738574040214
$ws.Range("A43").Value = 99.268738574040214
$ws.Range("A44").Value = 99.314442413162709
$ws.Range("A45").Value = 99.268738574040214
$ws.Range("A46").Value = 99.268738574040214
$ws.Range("A47").Value = 99.268738574040214
$ws.Range("A48").Value = 99.268738574040214
$ws.Range("A49").Value = 99.268738574040214
